$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the old layout completely (old table used A1:H3, incl. the now-removed
# "pagTermine" column and the separated "Action" column H) --------------------
$ws.Range("A1:H10").Clear()

# --- Apply the "Text" (@) number format to the whole working area first so that
# values which look like dates ("01.01.2026", "31.12.2025", ...) are stored as
# literal text instead of being auto-converted to date serials. -----------------
$ws.Range("A1:G4").NumberFormat = "@"

# --- Row 1: headers -------------------------------------------------------------
$ws.Range("A1").Value = "Control"
$ws.Range("B1").Value = "Modus"
$ws.Range("C1").Value = "calVon"
$ws.Range("D1").Value = "calBis"
$ws.Range("E1").Value = "mltKommentar"
$ws.Range("F1").Value = "butSpeichern"
$ws.Range("G1").Value = "Action"

# --- Row 2: Record/Selector ------------------------------------------------------
$ws.Range("A2").Value = "Record/Selector"
$ws.Range("B2").Value = "Modus"
$ws.Range("C2").Value = "id=from_date"
$ws.Range("D2").Value = "id=till_date"
$ws.Range("E2").Value = "id=comment"
$ws.Range("F2").Value = "xpath=//button[text()='Speichern']"
$ws.Range("G2").Value = "Action"

# Highlight the freshly recorded selectors in green.
$ws.Range("C2:F2").Interior.Color = 5296274

# --- Row 3: Check defaults -------------------------------------------------------
$ws.Range("A3").Value = "Check defaults"
$ws.Range("B3").Value = "<CHK>"
$ws.Range("C3").Value = "<EMPTY>"
$ws.Range("D3").Value = "<EMPTY>"
$ws.Range("E3").Value = "<EMPTY>"
$ws.Range("F3").Value = "<ENABLED>"
$ws.Range("G3").Value = "<NOP>"

# --- Row 4: Abwesenheit Sylvester anlegen (new test case row) -------------------
$ws.Range("A4").Value = "Abwesenheit Sylvester anlegen"
$ws.Range("B4").Value = "<SET>"
$ws.Range("C4").Value = "31.12.2025"
$ws.Range("D4").Value = "01.01.2026"
$ws.Range("E4").Value = "Hoch die Tassen!"
$ws.Range("F4").Value = "X"
$ws.Range("G4").Value = "<NOP>"

# --- Column widths ---------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 42.6640625
$ws.Columns.Item(2).ColumnWidth = 12.21875
$ws.Columns.Item(3).ColumnWidth = 26.21875
$ws.Columns.Item(4).ColumnWidth = 29.6640625
$ws.Columns.Item(5).ColumnWidth = 32.21875
$ws.Columns.Item(6).ColumnWidth = 29.33203125
$ws.Columns.Item(7).ColumnWidth = 11.5546875

# --- Sheet selection ---------------------------------------------------------------
$ws.Range("D6").Select()

# --- Page setup ---------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PrintQuality = 300

# --- Reposition the instructional screenshot -----------------------------------------
$shp = $ws.Shapes.Item(1)
$shp.Top = 92.4
$shp.Left = 0.6
$shp.Width = 1268.091496062992
$shp.Height = 569.1788188976378
